$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 178 (date 2021-02-11, K/L/M/P = 1500).
# This pushes the existing rows 178-180 down to 179-181, copying their
# formatting (incl. the date number-format on column D) along the way.
$ws.Rows.Item(178).Insert()

# Populate the freshly inserted row 178 with the new weekly price entry.
$ws.Cells.Item(178, 1).Value = 4
$ws.Cells.Item(178, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(178, 3).Value = "Los Lagos"
$ws.Cells.Item(178, 4).Value = 44448
$ws.Cells.Item(178, 5).Value = 10
$ws.Cells.Item(178, 6).Value = 100112008
$ws.Cells.Item(178, 7).Value = "Coliflor"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 500
$ws.Cells.Item(178, 11).Value = 1200
$ws.Cells.Item(178, 12).Value = 1200
$ws.Cells.Item(178, 13).Value = 1200
$ws.Cells.Item(178, 14).Value = "`$/unidad"
$ws.Cells.Item(178, 15).Value = "Región Metropolitana"
$ws.Cells.Item(178, 16).Value = 1200
$ws.Cells.Item(178, 17).Value = 1
$ws.Cells.Item(178, 18).Value = "Hortaliza"
